$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59 (shifts existing rows 59.. down by one)
$ws.Rows("59").Insert()

# Populate the newly inserted row 59 with its data
$ws.Range("A59").Value = 4
$ws.Range("B59").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C59").Value = "Los Lagos"
$ws.Range("D59").Value = 44495
$ws.Range("E59").Value = 10
$ws.Range("F59").Value = 100112017
$ws.Range("G59").Value = "Apio"
$ws.Range("H59").Value = "Americana (o)"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 40
$ws.Range("K59").Value = 10000
$ws.Range("L59").Value = 10000
$ws.Range("M59").Value = 10000
$ws.Range("N59").Value = "`$/docena de matas"
$ws.Range("O59").Value = "Región de Coquimbo"
$ws.Range("P59").Value = 1667
$ws.Range("Q59").Value = 6
$ws.Range("R59").Value = "Hortaliza"
